$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 36 becomes a new data row. Column A re-uses the existing "Diff"
# category label (already used by rows 30/32), column B/C are new entries.
$ws.Range("A36").Value = "Diff"
$ws.Range("B36").Value = "Extend git diff for excel (2)"
$ws.Range("C36").Value = "Previously we use python to parse an excel as git-diff's plugin. And now I try to use the liba.exe as the plugin and it works! Only need to change the textconv (in file .git/config) from previous setting:`n# textconv = ./liba -show"

# Match the row height used by the new wrapped content
$ws.Rows.Item(36).RowHeight = 39

# Move the active selection from C36 to C37, as recorded in the sheet view
$ws.Range("C37").Select()
